# Updates Price (D) and Volume(1h) (E) columns in the crypto symbol list
# to the latest scraped values, per the "Updated symbol list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.49%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "34.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "12.09%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.174"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "5.20%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07779"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "6.31%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.292"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.52%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.002"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.49%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.993"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "7.29%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9274"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.74%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "10.37%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1811"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "7.00%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08562"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.43%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03478"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "11.07%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09907"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.33%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001506"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.64%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04616"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.65%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005816"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.95%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.468"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.76%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.105"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.38%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "3.33%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1327"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.11%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.561"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "8.67%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2341"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "11.44%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001223"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.07%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004431"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.51%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.21%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "0.69%"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "11.43%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04719"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "6.11%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007594"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.43%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1407"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "6.08%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007066"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-25.47%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002215"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.23%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009195"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.22%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005990"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.08%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.20%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "6.71%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "35.10%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.20%"
